# Data_analysis. Calculo de acciones en USD
# Updates the "edn" quotes sheet:
#   - corrects montoOperado/volumenNominal/cantidadOperaciones for 2022-05-17 (row 93)
#   - corrects the 2024-05-27 row (row 583), whose trading stats had been
#     recorded incomplete
#   - appends five new daily rows (2024-05-29 .. 2024-06-04)
#
# NOTE: this interpreter's function calls only bind PARAMETERS POSITIONALLY
# (named `-Param value` binding does not work), so Set-Quote below is called
# positionally everywhere.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-Quote {
    param($Row, $Fecha, $Ultimo, $Apertura, $Maximo, $Minimo, $Monto, $Volumen, $Operaciones)

    if ($Fecha) {
        # Leading apostrophe forces the text-entry form Excel uses for
        # "looks like a date" strings, so the cell stays plain text instead
        # of being auto-converted to a date serial number.
        $ws.Cells.Item($Row, 1).Value = "'" + $Fecha
    }
    $ws.Cells.Item($Row, 2).Value2 = $Ultimo
    $ws.Cells.Item($Row, 3).Value2 = $Apertura
    $ws.Cells.Item($Row, 4).Value2 = $Maximo
    $ws.Cells.Item($Row, 5).Value2 = $Minimo
    $ws.Cells.Item($Row, 6).Value2 = $Monto
    $ws.Cells.Item($Row, 7).Value2 = $Volumen
    $ws.Cells.Item($Row, 8).Value2 = $Operaciones
}

# Row 93 (2022-05-17): only montoOperado / volumenNominal / cantidadOperaciones change
$apertura93 = $ws.Cells.Item(93, 3).Value2
$maximo93   = $ws.Cells.Item(93, 4).Value2
$minimo93   = $ws.Cells.Item(93, 5).Value2
$ultimo93   = $ws.Cells.Item(93, 2).Value2
Set-Quote 93 $null $ultimo93 $apertura93 $maximo93 $minimo93 1275578.1 35 105

# Row 583 (2024-05-27): ultimoPrecio / maximo / montoOperado / volumenNominal /
# cantidadOperaciones change; apertura / minimo stay as-is
$apertura583 = $ws.Cells.Item(583, 3).Value2
$minimo583   = $ws.Cells.Item(583, 5).Value2
Set-Quote 583 $null 1034 $apertura583 1036.9 $minimo583 289532561.7 283055 2078

# New rows appended after the previous last row (584)
Set-Quote 585 "2024-05-29" 1109.9 1100.5 1111.3  1077.3  665269844.15 605148 3017
Set-Quote 586 "2024-05-30" 1135.5 1115   1165    1106    825385513.45 728997 3127
Set-Quote 587 "2024-05-31" 1109   1139   1139    1093    763165653.8  685227 2839
Set-Quote 588 "2024-06-03" 1128   1150   1150    1095.8  293258102.5  1      1257
Set-Quote 589 "2024-06-04" 1070   1105   1108.75 1059    454474081.55 423284 2494
